$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number => ordered column/value pairs reflecting the
# corrected match data (id, teams, score, result, odds) for that row.
$rowData = @{}

$rowData[17] = [ordered]@{"B"=6221640; "F"="Ordabasy"; "G"="Zhetysu"; "H"=2; "I"=0; "J"="H"; "K"=1.333; "L"=4.333; "M"=8; "N"=1.333; "O"=4.333; "P"=8; "Q"=-1.25; "R"=1.75; "S"=1.95; "T"=2.5; "U"=1.9; "V"=1.9; "W"=0.333; "X"=-1; "Y"=-1; "Z"=0.75; "AA"=-1; "AB"=-1; "AC"=0.8999999999999999}
$rowData[18] = [ordered]@{"B"=6221766; "F"="Kairat Almaty"; "G"="FK Kaspyi Aktau"; "H"=3; "I"=1; "J"="H"; "K"=1.55; "L"=3.8; "M"=5; "N"=1.65; "O"=4; "P"=4.5; "Q"=-0.75; "R"=1.8; "S"=2; "T"=2.75; "U"=1.925; "V"=1.875; "W"=0.6499999999999999; "X"=-1; "Y"=-1; "Z"=0.8; "AA"=-1; "AB"=0.925; "AC"=-1}
$rowData[19] = [ordered]@{"B"=6221639; "F"="FC Astana"; "G"="FK Aktobe"; "H"=1; "I"=4; "J"="A"; "K"=1.45; "L"=4; "M"=5.75; "N"=1.533; "O"=3.75; "P"=5; "Q"=-1; "R"=1.9; "S"=1.9; "T"=2.75; "U"=1.95; "V"=1.85; "W"=-1; "X"=-1; "Y"=4; "Z"=-1; "AA"=0.8999999999999999; "AB"=0.95; "AC"=-1}
$rowData[20] = [ordered]@{"B"=6221642; "F"="Kaisar Kyzylorda"; "G"="FK Aksu"; "H"=2; "I"=2; "J"="D"; "K"=2.45; "L"=3; "M"=2.7; "N"=2.15; "O"=3.3; "P"=3.1; "Q"=-0.25; "R"=1.9; "S"=1.9; "T"=2.5; "U"=1.925; "V"=1.875; "W"=-1; "X"=2.3; "Y"=-1; "Z"=-0.5; "AA"=0.45; "AB"=0.925; "AC"=-1}
$rowData[27] = [ordered]@{"B"=6221769; "F"="FK Atyrau"; "G"="Okzhetpes Kokshetau"; "H"=0; "I"=4; "J"="A"; "K"=1.75; "L"=3.4; "M"=4; "N"=1.75; "O"=3.4; "P"=4; "Q"=-0.5; "R"=1.8; "S"=2; "T"=2.5; "U"=1.95; "V"=1.85; "W"=-1; "X"=-1; "Y"=3; "Z"=-1; "AA"=1; "AB"=0.95; "AC"=-1}
$rowData[28] = [ordered]@{"B"=6221645; "F"="FK Aktobe"; "G"="FK Kyzylzhar"; "H"=1; "I"=3; "J"="A"; "K"=1.5; "L"=4; "M"=5; "N"=1.727; "O"=3.5; "P"=4; "Q"=-0.75; "R"=2; "S"=1.8; "T"=2.25; "U"=1.85; "V"=1.95; "W"=-1; "X"=-1; "Y"=3; "Z"=-1; "AA"=0.8; "AB"=0.8500000000000001; "AC"=-1}
$rowData[38] = [ordered]@{"B"=6221774; "F"="Shakhter Karagandy"; "G"="FK Maktaaral"; "H"=3; "I"=0; "J"="H"; "K"=2.6; "L"=3; "M"=2.5; "N"=2.7; "O"=3; "P"=2.4; "Q"=0; "R"=2.025; "S"=1.775; "T"=2.25; "U"=1.85; "V"=1.95; "W"=1.7; "X"=-1; "Y"=-1; "Z"=1.025; "AA"=-1; "AB"=0.8500000000000001; "AC"=-1}
$rowData[39] = [ordered]@{"B"=6221655; "F"="Zhetysu"; "G"="Kaisar Kyzylorda"; "H"=0; "I"=2; "J"="A"; "K"=2.4; "L"=3.1; "M"=2.7; "N"=1.95; "O"=3.4; "P"=3.4; "Q"=-0.25; "R"=1.75; "S"=1.95; "T"=2.25; "U"=1.8; "V"=2; "W"=-1; "X"=-1; "Y"=2.4; "Z"=-1; "AA"=0.95; "AB"=-0.5; "AC"=0.5}
$rowData[177] = [ordered]@{"B"=6221753; "F"="FK Aksu"; "G"="Tobol Kostanay"; "H"=0; "I"=3; "J"="A"; "K"=2.75; "L"=3.1; "M"=2.375; "N"=2.625; "O"=3.2; "P"=2.45; "Q"=0; "R"=2; "S"=1.8; "T"=2.5; "U"=1.9; "V"=1.9; "W"=-1; "X"=-1; "Y"=1.45; "Z"=-1; "AA"=0.8; "AB"=0.8999999999999999; "AC"=-1}
$rowData[178] = [ordered]@{"B"=6221815; "F"="FK Atyrau"; "G"="Kairat Almaty"; "H"=0; "I"=0; "J"="D"; "K"=3; "L"=3; "M"=2.25; "N"=3.1; "O"=3.1; "P"=2.15; "Q"=0.25; "R"=1.85; "S"=1.95; "T"=2.25; "U"=1.8; "V"=2; "W"=-1; "X"=2.1; "Y"=-1; "Z"=0.425; "AA"=-0.5; "AB"=-1; "AC"=1}
$rowData[179] = [ordered]@{"B"=6221752; "F"="FK Kyzylzhar"; "G"="Kaisar Kyzylorda"; "H"=0; "I"=1; "J"="A"; "K"=1.833; "L"=3.2; "M"=4; "N"=1.85; "O"=3.2; "P"=4; "Q"=-0.5; "R"=1.9; "S"=1.9; "T"=2; "U"=1.775; "V"=2.025; "W"=-1; "X"=-1; "Y"=3; "Z"=-1; "AA"=0.8999999999999999; "AB"=-1; "AC"=-1}

foreach ($rowNum in $rowData.Keys) {
    $cols = $rowData[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $colIndex = $ws.Range($colLetter + "1").Column
        $ws.Cells.Item($rowNum, $colIndex).Value = $cols[$colLetter]
    }
}

"Updated " + $rowData.Keys.Count + " rows"
